# Creating send notifications file (Buyers and Emails)
#
# Repurposes the old "Notifications By Buyers" input-file block (rows 14-16
# on the Settings sheet) into a new "Send Notifications" input-file block,
# adds a new row for the Emails sheet name (row 17), and adds a brand new
# "Send SMTP Mail Message" settings block (rows 20-23).
# Also flips which sheet is active/selected (Assets becomes the active tab
# instead of Settings).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws3 = $wb.Worksheets.Item("Assets")

# --- Rows 14-17: rename the "Notifications by Buyers" input file block into
#     the new "Send Notifications" input file block. Values are written in
#     this particular order so the shared-string table is rebuilt with the
#     same layout/index assignment Excel itself produced.
$ws1.Range("A14").Value = "Input_SendNotifications (Input File)"

$ws1.Range("A15").Value = "File_SendNotificationsName"
$ws1.Range("C15").Value = "Input file - Notifications by Buyers"

$ws1.Range("A16").Value = "Sheet_BuyersName"
$ws1.Range("B16").Value = "Buyers"

$ws1.Range("A17").Value = "Sheet_EmailsName"
$ws1.Range("B17").Value = "Emails"
$ws1.Range("C17").Value = "Sheet name for Send Notifications"

$ws1.Range("C16").Value = "Sheet name for Send Notifications by buyers"
$ws1.Range("B15").Value = "Input_SendNotifications.xlsx"

# --- Row 20: new section header, styled like the other section headers
#     (e.g. row 8) by copying that formatting across.
$ws1.Range("A8:C8").Copy()
$ws1.Range("A20:C20").PasteSpecial(-4122)
$ws1.Range("A20").Value = "Send SMTP Mail Message"

# --- Rows 21-23: new SMTP settings (port / host / account).
$ws1.Range("A21").Value = "EmailSMTPServerPortName"
$ws1.Range("B21").HorizontalAlignment = -4131
$ws1.Range("B21").VerticalAlignment = -4108
$ws1.Range("B21").Value = 25

$ws1.Range("B21").Copy()
$ws1.Range("B22").PasteSpecial(-4122)
$ws1.Range("B22").Value = "10.101.1.126"
$ws1.Range("A22").Value = "EmailSMTPServerHostName"

$ws1.Range("A23").Value = "EmailSMTPAccountName"

# --- Switch the active sheet/selection: Assets becomes the selected tab
#     (previously Settings was), with a fresh selection on each sheet.
$ws1.Range("B23").Select()
$ws3.Activate()
$ws3.Range("A5").Select()
